$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 413-414; this shifts the existing rows 413-495
# down to 415-497 (formats/number formats are inherited from the row
# above, matching the date style already used in column D).
$ws.Rows("413:414").Insert()

# Row 413: new Murcott "Primera" price entry dated 45258 (2023-11-28)
$ws.Range("A413").Value = 7
$ws.Range("B413").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C413").Value = "Ñuble"
$ws.Range("D413").Value2 = 45258
$ws.Range("E413").Value = 16
$ws.Range("F413").Value = "Fruta"
$ws.Range("G413").Value = 100102
$ws.Range("H413").Value = "Cítricos"
$ws.Range("I413").Value = 100102004
$ws.Range("J413").Value = "Mandarina"
$ws.Range("K413").Value = "Murcott"
$ws.Range("L413").Value = "Primera"
$ws.Range("M413").Value = 120
$ws.Range("N413").Value = 9000
$ws.Range("O413").Value = 10000
$ws.Range("P413").Value = 9500
$ws.Range("Q413").Value = "$/bandeja 18 kilos"
$ws.Range("R413").Value = "Región de O'Higgins"
$ws.Range("S413").Value = 528
$ws.Range("T413").Value = 18

# Row 414: new Murcott "Segunda" price entry, same date
$ws.Range("A414").Value = 7
$ws.Range("B414").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C414").Value = "Ñuble"
$ws.Range("D414").Value2 = 45258
$ws.Range("E414").Value = 16
$ws.Range("F414").Value = "Fruta"
$ws.Range("G414").Value = 100102
$ws.Range("H414").Value = "Cítricos"
$ws.Range("I414").Value = 100102004
$ws.Range("J414").Value = "Mandarina"
$ws.Range("K414").Value = "Murcott"
$ws.Range("L414").Value = "Segunda"
$ws.Range("M414").Value = 100
$ws.Range("N414").Value = 8000
$ws.Range("O414").Value = 8000
$ws.Range("P414").Value = 8000
$ws.Range("Q414").Value = "$/bandeja 18 kilos"
$ws.Range("R414").Value = "Región de O'Higgins"
$ws.Range("S414").Value = 444
$ws.Range("T414").Value = 18
